$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRF = 11.02637931034483

for ($r = 34; $r -le 52; $r++) {
    $ws.Cells.Item($r, 9).Value = $newRF
}
